# working_hours.xlsx: one more pair of coefficients (a new shift) was
# measured for 2014-03-27/28, so two more data rows are inserted right
# before the summary block. The summary rows (sum [min] / sum [h] /
# sum [working weeks]) shift down by two rows and their formulas - which
# reference the data range - grow to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last data row was 97 (2014-03-27), followed by a blank
# spacer row (98) and the three summary rows (99-101). Insert two fresh
# rows right after the last populated data row so the existing spacer /
# summary rows (and their formulas, e.g. "=SUM(F2:F98)") shift down and
# their ranges grow automatically, exactly like typing new rows in the
# middle of the table in Excel.
$ws.Rows.Item(98).Insert() | Out-Null
$ws.Rows.Item(99).Insert() | Out-Null

# New data row: 2014-03-27, 16:50 -> 19:15
$ws.Cells.Item(98, 1).Value = 2014
$ws.Cells.Item(98, 2).Value = 3
$ws.Cells.Item(98, 3).Value = 27
$ws.Cells.Item(98, 4).Value = 0.70138888888888884
$ws.Cells.Item(98, 5).Value = 0.80208333333333337
$ws.Cells.Item(98, 6).Formula = "=(E98-D98)*24*60"
$ws.Cells.Item(98, 7).Formula = "=F98/60"

# New data row: 2014-03-28, 08:20 -> 09:45
$ws.Cells.Item(99, 1).Value = 2014
$ws.Cells.Item(99, 2).Value = 3
$ws.Cells.Item(99, 3).Value = 28
$ws.Cells.Item(99, 4).Value = 0.34722222222222227
$ws.Cells.Item(99, 5).Value = 0.40625
$ws.Cells.Item(99, 6).Formula = "=(E99-D99)*24*60"
$ws.Cells.Item(99, 7).Formula = "=F99/60"

# Match the saved selection from the edit (active cell moved to the new
# "sum [min]" total cell one column/row further down).
[void]$ws.Range("F99").Select()

$wb.Application.Calculate() | Out-Null
